$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: Validation -> F1 train
$ws.Range("O1").Value = "F1 train"

# Row 2
$ws.Range("O2").Value = 1

# Row 3
$ws.Range("O3").Value = 1

# Row 4
$ws.Range("O4").Value = 0.958904109589041

# Row 5
$ws.Range("O5").Value = 1

# Row 6 (MLP params + metrics changed)
$ws.Range("C6").Value = "{'activation': 'tanh', 'alpha': 0.0001, 'hidden_layer_sizes': (64,), 'learning_rate': 'constant'}"
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 0.7
$ws.Range("J6").Value = 0.75
$ws.Range("K6").Value = 0.9
$ws.Range("L6").Value = 0.6428571428571429
$ws.Range("M6").Value = 0.5
$ws.Range("N6").Value = 0.9
$ws.Range("O6").Value = 0.6024096385542169

# Row 7
$ws.Range("O7").Value = 0.9736842105263158

# Row 8
$ws.Range("O8").Value = 1

# Row 9
$ws.Range("O9").Value = 0.9210526315789473

# Row 10
$ws.Range("O10").Value = 1

# Row 11
$ws.Range("O11").Value = 0.5818181818181818

# Row 12
$ws.Range("O12").Value = 1

# Row 13
$ws.Range("O13").Value = 1

# Row 14
$ws.Range("O14").Value = 1

# Row 15
$ws.Range("O15").Value = 0.7605633802816901

# Row 16 (MLP params + metrics changed)
$ws.Range("C16").Value = "{'activation': 'relu', 'alpha': 0.0001, 'hidden_layer_sizes': (64, 32), 'learning_rate': 'constant'}"
$ws.Range("E16").Value = 5
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 5
$ws.Range("I16").Value = 0.75
$ws.Range("K16").Value = 0.5
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1
$ws.Range("N16").Value = 0.5
$ws.Range("O16").Value = 0.5
